$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes old rows 2-4 down to 3-5)
$ws.Rows("2:2").Insert()

# --- Header row (row 1): rename J1, add new headers K1 and L1 ---
$ws.Range("J1").Value = "eAniSet"
$ws.Range("K1").Value = "minPatX"
$ws.Range("L1").Value = "maxPatX"

# --- Row 2 (new spawner entry) ---
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 40
$ws.Range("C2").Value = 35
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 340
$ws.Range("F2").Value = 150
$ws.Range("G2").Value = 20
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 134
$ws.Range("J2").Value = 88
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 330

# --- Row 3 (previously row 2): strip formulas, keep the shifted values ---
$ws.Range("B3").Value = 368
$ws.Range("C3").Value = 34
$ws.Range("E3").Value = 688
$ws.Range("F3").Value = 150
$ws.Range("K3").Value = 390
$ws.Range("L3").Value = 600

# --- Row 4 (previously row 3): strip formulas, keep the shifted values ---
$ws.Range("B4").Value = 672
$ws.Range("C4").Value = 15
$ws.Range("E4").Value = 992
$ws.Range("F4").Value = 131
$ws.Range("K4").Value = 680
$ws.Range("L4").Value = 1080

# --- Row 5 (previously row 4): strip formulas, keep the shifted values ---
$ws.Range("B5").Value = 720
$ws.Range("C5").Value = 15
$ws.Range("E5").Value = 1040
$ws.Range("F5").Value = 131
$ws.Range("K5").Value = 680
$ws.Range("L5").Value = 1080

# --- Update selection / view to match the saved workbook state ---
$ws.Range("A2:L5").Select()

$wb.Save()
